$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.5772774815559387
$ws.Range("B1").Value = 1.183201789855957
$ws.Range("C1").Value = 15
$ws.Range("D1").Value = 2.733784675598145
$ws.Range("E1").Value = 1.452551245689392
